# Updates the cryptos list data per the Oct 27 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) and Volume(1h) (E) columns hold plain-text values (including
# values that look numeric, e.g. "589.93" or "2.512.40"). Force the range
# to a Text number format before writing so Excel doesn't silently coerce
# them into floating point numbers, then restore the default "Normal"
# style so no stray formatting is left behind on the cells.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "68.076.70"
$ws.Range("E2").Value = "  +1.44%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.512.40"
$ws.Range("E3").Value = "  +1.17%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.07%  "

# Row 5 - BNB
$ws.Range("D5").Value = "589.93"
$ws.Range("E5").Value = "  +1.07%  "

# Row 6 - Solana
$ws.Range("D6").Value = "177.22"
$ws.Range("E6").Value = "  +3.65%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.07%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  +0.81%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.145"
$ws.Range("E9").Value = "  +5.17%  "

# Row 10 - TRON
$ws.Range("E10").Value = "  -0.56%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +2.07%  "

# Row 12 - Toncoin
$ws.Range("E12").Value = "  +0.49%  "

# Row 14 - Avalanche
$ws.Range("D14").Value = "25.78"
$ws.Range("E14").Value = "  +1.48%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "68.008.94"
$ws.Range("E15").Value = "  +1.62%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  +1.28%  "

# Row 17 - WrappedEther
$ws.Range("D17").Value = "2.495.38"
$ws.Range("E17").Value = "  +2.87%  "

# Row 18 - Chainlink
$ws.Range("D18").Value = "11.04"
$ws.Range("E18").Value = "  +0.22%  "

# Row 19 - Uniswap
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +1.45%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "352.80"
$ws.Range("E20").Value = "  +1.11%  "

# Row 21 - Polkadot
$ws.Range("E21").Value = "  +2.40%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.01%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "70.99"
$ws.Range("E23").Value = "  +3.80%  "

# Row 24 - NEARProtocol
$ws.Range("D24").Value = "4.29"
$ws.Range("E24").Value = "  +1.27%  "

# Row 25 - SuiNetwork
$ws.Range("D25").Value = "1.76"
$ws.Range("E25").Value = "  -2.34%  "

# Row 26 - Aptos
$ws.Range("D26").Value = "9.19"
$ws.Range("E26").Value = "  -1.21%  "

# Row 27 - WrappedeETH
$ws.Range("D27").Value = "2.639.36"
$ws.Range("E27").Value = "  +0.67%  "

# Row 28 - Binance-PegBSC-USD
$ws.Range("D28").Value = "0.986"
$ws.Range("E28").Value = "  -1.44%  "

# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0921"
$ws.Range("E29").Value = "  +1.62%  "

# Row 30 - Bittensor
$ws.Range("D30").Value = "508.94"
$ws.Range("E30").Value = "  -0.76%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "7.82"
$ws.Range("E31").Value = "  +1.48%  "

# Row 32 - Fetch.AI
$ws.Range("D32").Value = "1.27"
$ws.Range("E32").Value = "  +2.94%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +1.06%  "

# Row 34 - FirstDigitalUSD
$ws.Range("E34").Value = "  -0.03%  "

# Row 35 - Kaspa
$ws.Range("E35").Value = "  +3.48%  "

# Row 36 - Monero
$ws.Range("D36").Value = "164.61"
$ws.Range("E36").Value = "  +2.51%  "

# Row 37 - EthereumClassic
$ws.Range("D37").Value = "18.44"
$ws.Range("E37").Value = "  +1.21%  "

# Row 38 - WhiteBITCoin
$ws.Range("D38").Value = "18.65"
$ws.Range("E38").Value = "  -0.32%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  +0.49%  "

# Row 40 - USDe
$ws.Range("E40").Value = "  +0.04%  "

# Row 41 - Stacks
$ws.Range("E41").Value = "  +3.38%  "

# Row 42 - was RenderToken, now PolygonEcosystemToken
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").Value = "0.331"
$ws.Range("E42").Value = "  +0.63%  "

# Row 43 - was PolygonEcosystemToken, now RenderToken
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.90"
$ws.Range("E43").Value = "  +2.07%  "

# Row 44 - dogwifhat
$ws.Range("D44").Value = "2.50"
$ws.Range("E44").Value = "  +5.77%  "

# Row 45 - Aave
$ws.Range("D45").Value = "147.87"
$ws.Range("E45").Value = "  +3.51%  "

# Row 46 - Filecoin
$ws.Range("D46").Value = "3.57"
$ws.Range("E46").Value = "  +3.12%  "

# Row 47 - was ARBITRUM, now BabyDogeCoin
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0263"
$ws.Range("E47").Value = "  +4.84%  "

# Row 48 - was BabyDogeCoin, now ARBITRUM
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").Value = "0.521"
$ws.Range("E48").Value = "  +1.28%  "

# Row 49 - Cronos
$ws.Range("E49").Value = "  +2.30%  "

# Row 50 - Optimism
$ws.Range("E50").Value = "  +2.15%  "

# Row 51 - Mantle
$ws.Range("D51").Value = "0.588"
$ws.Range("E51").Value = "  +0.96%  "

# Restore the default cell style so no stray number-format styling is
# left on the (now textual) Price/Volume cells.
$dataRange.Style = "Normal"
